# Weekly update: a new price record is inserted before the current row 27,
# pushing the existing rows 27-49 down to 28-50 (dimension grows to A1:T50).
# The newly inserted row 27 receives a fresh data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 27; rows 27-49 shift down to 28-50.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Range("A27").Value = 6
$ws.Range("B27").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C27").Value = 'Metropolitana'
$ws.Range("D27").Value = 45033
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100101
$ws.Range("H27").Value = 'Berries'
$ws.Range("I27").Value = 100101006
$ws.Range("J27").Value = 'Higo'
$ws.Range("K27").Value = 'Sin especificar'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 20000
$ws.Range("Q27").Value = '$/bandeja 7 kilos'
$ws.Range("R27").Value = 'Región Metropolitana'
$ws.Range("S27").Value = 2857
$ws.Range("T27").Value = 7
